$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$shp = $ws.Shapes.Item(1)
Write-Host "before Left:" $shp.Left
$shp.Left = $shp.Left + 7.5
Write-Host "after Left:" $shp.Left
